# Apply updated FedEx tracking numbers / status (DEV URL configuration change)
# Re-run of the shipment export produced new tracking numbers (col P) for every
# row, and row 20's shipment now reports its correct dollar amount (col Q) and
# flips from FAIL to PASS (col R).
#
# Column P holds numeric-looking tracking numbers that must stay TEXT (shared
# string), matching the source file. Typing a pure digit string into a
# General-formatted cell makes Excel coerce it to a real number, so we route
# the value through a helper cell holding a text-producing formula and paste
# only the (text) VALUE into the destination - this preserves the
# destination's existing (unstyled) cell format exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($targetAddr, $text) {
    $scratch = $ws.Range("D2")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
    $scratch.Value = ""
}

Set-TextValue "P2"  "320018799382"
Set-TextValue "P3"  "320018799393"
Set-TextValue "P4"  "320018799420"
Set-TextValue "P5"  "320018799441"
Set-TextValue "P6"  "320018799485"
Set-TextValue "P7"  "320018799500"
Set-TextValue "P8"  "320018799533"
Set-TextValue "P9"  "320018799625"
Set-TextValue "P10" "320018799658"
Set-TextValue "P11" "320018799670"
Set-TextValue "P12" "320018799717"
Set-TextValue "P13" "320018799739"
Set-TextValue "P14" "320018799761"
Set-TextValue "P15" "320018799783"
Set-TextValue "P16" "320018799810"
Set-TextValue "P17" "320018799831"
Set-TextValue "P18" "320018799875"
Set-TextValue "P19" "320018792701"
Set-TextValue "P20" "320018792734"
Set-TextValue "P21" "320018792756"
Set-TextValue "P22" "320018792789"
Set-TextValue "P23" "320018792790"
Set-TextValue "P24" "320018792804"
Set-TextValue "P25" "320018792815"
Set-TextValue "P26" "320018792826"

# Row 20: amount re-priced and validation flipped from FAIL to PASS
Set-TextValue "Q20" "$62.39"
$ws.Range("R20").Value = "PASS"
